# 7.8 History Card & Advanced Story
# Updates the "Paper" investigation dialogue table on Sheet1:
#  - inserts a new row (Dee: "Let's investigate more closely.") before the
#    old "Investigate/Desk" row
#  - refreshes several lines of dialogue text
#  - adds Avatar (column C) / BGM (column F) cues that were missing
#  - renames the "Lee-Regular" cue to "Dee-Thinking" and retargets the
#    "Master" reference to "Lord"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 7 - this pushes the old rows 7-11 (Investigate/Desk,
# Paper, Book, End Investigation, disappear) down to rows 8-12 and keeps
# their existing values/styles intact.
$ws.Rows.Item(7).EntireRow.Insert()
$ws.Rows.Item(7).RowHeight = 17

# Row 2 - Yao: "Sir, there are several stacks ..." (add Avatar + BGM cues,
# rename the Character2 cue, keep the History/appearAt/500 trigger)
$ws.Cells.Item(2, 3).Value = "Yao-Regular"
$ws.Cells.Item(2, 6).Value = "Suspicious"
$ws.Cells.Item(2, 12).Value = "Dee-Thinking"

# Row 3 - Yao: "But this paper looks quite different ..." (add Avatar cue)
$ws.Cells.Item(3, 3).Value = "Yao-Query"

# Row 4 - He: custom-made paper line (add Avatar cue, tweak punctuation)
$ws.Cells.Item(4, 2).Value = "Indeed, this is custom-made paper used exclusively at Qingliu Manor——it’s not available outside."
$ws.Cells.Item(4, 3).Value = "He-Regular"

# Row 5 - He: "Even within the manor ..." (Master -> Lord, add Avatar cue)
$ws.Cells.Item(5, 2).Value = "Even within the manor, only the Lord had access to it."
$ws.Cells.Item(5, 3).Value = "He-Regular"

# Row 6 - Dee: shorten to just the "I see! ..." line, add Avatar cue
$ws.Cells.Item(6, 2).Value = "I see! I feel like I’ve seen this type of paper somewhere before......"
$ws.Cells.Item(6, 3).Value = "Dee-Determined"

# Row 7 (new) - Dee: the "Let's investigate more closely." line that used
# to be appended to row 6
$ws.Cells.Item(7, 1).Value = "Dee"
$ws.Cells.Item(7, 2).Value = "Let’s investigate more closely."
$ws.Cells.Item(7, 3).Value = "Dee-Determined"
$ws.Cells.Item(7, 4).Value = "DialogueVocal"
$ws.Cells.Item(7, 5).Value = "StudyInvestigate"

# Tidy up the selection/used-range bookkeeping to match the new row count
$ws.Range("B13").Select()
